$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Blocco_1')
$ws.Range('C57').Value = "'531877"
$ws.Range('D57').Value = "'C.so Francia   10093 Collegno ( TO )"
$ws.Range('E57').Value = "'FIBE srl"
$ws.Range('C58').Value = "'538088"
$ws.Range('D58').Value = "'C.so Francia  10093 Collegno ( TO )"
$ws.Range('E58').Value = "'DueB Costruzioni srl"
$ws = $wb.Worksheets.Item('Blocco_2')
$ws.Range('C69').Value = "'534139"
$ws.Range('D69').Value = "'Via Callas  35020 Albignasego ( PD )"
$ws.Range('E69').Value = "'Edilbaraldo srl"
$ws.Range('C70').Value = "'534138"
$ws.Range('D70').Value = "'Via M. Callas  35020 Albignasego ( PD )"
$ws.Range('E70').Value = "'Immobiliare San Bonaventura srl; assegnato"
$ws = $wb.Worksheets.Item('Blocco_3')
$ws.Range('C41').Value = "'542609"
$ws.Range('D41').Value = "'V.le delle Americhe  48122 Ravenna ( RA )"
$ws.Range('C42').Value = "'537153"
$ws.Range('D42').Value = "'Via delle Americhe  48122 Ravenna ( RA )"
$ws.Range('C48').Value = "'539779"
$ws.Range('D48').Value = "'Via Ravenna ang. Via Savio   47814 Bellaria-Igea Marina ( RN )"
$ws.Range('E48').Value = "'Moma srl"
$ws.Range('C49').Value = "'532124"
$ws.Range('D49').Value = "'Via Savio snc  47814 Bellaria-Igea Marina ( RN )"
$ws.Range('E49').Value = "'Pompili Daniele"
$ws.Range('F50').Value = 86688
$ws.Range('F51').Value = 88131
$ws.Range('F52').Value = 89760
$ws.Range('F53').Value = 91260
$ws.Range('F54').Value = 92938
$ws.Range('F55').Value = 94386
$ws.Range('F56').Value = 95753
$ws.Range('F57').Value = 96974
$ws.Range('F58').Value = 98347
$ws.Range('F59').Value = 99862
$ws.Range('F60').Value = 101671
$ws.Range('F61').Value = 103632
$ws.Range('F62').Value = 105007
$ws.Range('F63').Value = 106491
$ws.Range('F64').Value = 108589
$ws.Range('F65').Value = 110133
$ws.Range('F66').Value = 111734
$ws.Range('F67').Value = 113720
$ws.Range('F68').Value = 115282
$ws.Range('F69').Value = 117049
$ws.Range('F70').Value = 118665
$ws.Range('F71').Value = 121519
$ws.Range('F72').Value = 125883
$ws.Range('F73').Value = 129206
$ws.Range('F74').Value = 129206
$ws.Range('F75').Value = 131523
$ws.Range('F76').Value = 134535
$ws.Range('F77').Value = 136359
$ws.Range('F78').Value = 137559
$ws.Range('F79').Value = 139224
$ws.Range('F80').Value = 141276
$ws.Range('F81').Value = 142476
$ws.Range('F82').Value = 143871
$ws.Range('F83').Value = 145739
$ws.Range('F84').Value = 148260
$ws.Range('F85').Value = 150563
$ws.Range('F86').Value = 152661
$ws.Range('F87').Value = 154343
$ws.Range('F88').Value = 156752
$ws.Range('F89').Value = 159148
$ws.Range('F90').Value = 160906
$ws.Range('F91').Value = 162583
$ws.Range('F92').Value = 164267
$ws.Range('F93').Value = 165542
$ws.Range('F94').Value = 167190
$ws.Range('F95').Value = 169306
$ws.Range('F96').Value = 171906
$ws.Range('F97').Value = 173837
$ws.Range('F98').Value = 175140
$ws.Range('F99').Value = 176514
$ws.Range('F100').Value = 177714
$ws.Range('F101').Value = 179143
$ws.Range('F102').Value = 180409
$ws.Range('F103').Value = 181923
$ws.Range('F104').Value = 183606
$ws.Range('F105').Value = 185114
$ws.Range('F106').Value = 186733
$ws.Range('F107').Value = 188665
$ws.Range('F108').Value = 190307
$ws.Range('F109').Value = 191699
$ws.Range('F110').Value = 193078
$ws.Range('F111').Value = 195373
$ws.Range('F112').Value = 197560
$ws.Range('F113').Value = 199871
$ws.Range('F114').Value = 201206
$ws.Range('F115').Value = 202927
$ws.Range('F116').Value = 204778
$ws.Range('F117').Value = 206173
$ws.Range('F118').Value = 207673
$ws.Range('F119').Value = 209506
$ws.Range('F120').Value = 212181
$ws.Range('F121').Value = 214311
$ws.Range('F122').Value = 216291
$ws.Range('F123').Value = 218171
$ws.Range('F124').Value = 220163
$ws.Range('F125').Value = 221993
$ws.Range('F126').Value = 223980
$ws.Range('F127').Value = 226742
$ws.Range('F128').Value = 229313
$ws.Range('F129').Value = 230991
$ws.Range('F130').Value = 232333
$ws.Range('F131').Value = 236803
$ws.Range('F132').Value = 239296
$ws.Range('F133').Value = 240496
$ws.Range('F134').Value = 244422
$ws.Range('F135').Value = 245691
$ws.Range('F136').Value = 248427
$ws.Range('F137').Value = 252049
$ws.Range('F138').Value = 253692
$ws.Range('F139').Value = 255218
$ws.Range('F140').Value = 256670
